$wb = $excel.ActiveWorkbook

# --- Rename "Sheet1" to "debit spread 100B-200S" ---
$wsDebit = $wb.Worksheets.Item("Sheet1")
$wsDebit.Name = "debit spread 100B-200S"

# --- Add the new label cell J2 referencing the new shared string ---
$wsDebit.Range("J2").Value = "Debit spread 100 buy - 200 sell"

# --- Fix up the chart on that sheet so its series formulas point at the renamed sheet ---
$chart = $wsDebit.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('debit spread 100B-200S'!`$E`$1,,'debit spread 100B-200S'!`$E`$2:`$E`$21,1)"

# --- Update selections / active sheet to match the saved UI state ---
$wsDist200 = $wb.Worksheets.Item("200 dist")
$wsDist200.Activate()
$null = $wsDist200.Range("V32").Select()

$wsDist100 = $wb.Worksheets.Item("100 dist")
$wsDist100.Activate()
$null = $wsDist100.Range("D34").Select()

$wsDebit.Activate()
$null = $wsDebit.Range("M26").Select()

$wsOutput = $wb.Worksheets.Item("output")
$wsOutput.Activate()
